$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "25.854.32"
$ws.Range("E2").Value = "  -2.36%  "

# Row 3
$ws.Range("D3").Value = "1.752.61"
$ws.Range("E3").Value = "  -4.59%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9987"
$ws.Range("E4").Value = "  -0.18%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.71"
$ws.Range("E5").Value = "  -8.23%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9990"
$ws.Range("E6").Value = "  -0.16%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5110"
$ws.Range("E7").Value = "  -4.99%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.21"
$ws.Range("E8").Value = "  -5.89%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2788"
$ws.Range("E9").Value = "  -6.40%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06209"
$ws.Range("E10").Value = "  -10.35%  "

# Row 11
$ws.Range("D11").Value = "1.748.19"
$ws.Range("E11").Value = "  -4.92%  "

# Row 12
$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.81"
$ws.Range("E12").Value = "  -9.77%  "

# Row 13
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.06970"
$ws.Range("E13").Value = "  -3.43%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6099"
$ws.Range("E14").Value = "  -16.76%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.523"
$ws.Range("E15").Value = "  -9.36%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "77.71"
$ws.Range("E16").Value = "  -12.74%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9978"
$ws.Range("E17").Value = "  -0.31%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9989"
$ws.Range("E18").Value = "  -0.17%  "

# Row 19
$ws.Range("D19").Value = "25.846.23"
$ws.Range("E19").Value = "  -2.47%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006993"
$ws.Range("E20").Value = "  -11.62%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.71"
$ws.Range("E21").Value = "  -15.27%  "

# Row 22
$ws.Range("D22").Value = "1.970.22"
$ws.Range("E22").Value = "  -5.25%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.088"
$ws.Range("E23").Value = "  -10.79%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.264"
$ws.Range("E24").Value = "  -12.14%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.227"

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "137.81"
$ws.Range("E26").Value = "  -3.25%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.472"
$ws.Range("E27").Value = "  -13.72%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.829"
$ws.Range("E28").Value = "  -15.77%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.07"
$ws.Range("E29").Value = "  -11.23%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "103.88"
$ws.Range("E30").Value = "  -6.27%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08192"
$ws.Range("E31").Value = "  -7.52%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.705"
$ws.Range("E32").Value = "  -12.47%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.534"
$ws.Range("E33").Value = "  -12.43%  "

# Row 34
$ws.Range("E34").Value = "  -6.34%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9978"
$ws.Range("E35").Value = "  -0.20%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.626"
$ws.Range("E36").Value = "  -10.10%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9946"
$ws.Range("E37").Value = "  -12.05%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6133"
$ws.Range("E38").Value = "  -15.38%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.696"
$ws.Range("E39").Value = "  -12.67%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01559"
$ws.Range("E40").Value = "  -8.94%  "

# Row 41
$ws.Range("E41").Value = "  -16.66%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9994"
$ws.Range("E42").Value = "  -0.09%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "103.45"
$ws.Range("E43").Value = "  -3.84%  "

# Row 44
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3886"
$ws.Range("E44").Value = "  -17.31%  "

# Row 45
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7489"
$ws.Range("E45").Value = "  -17.06%  "

# Row 46
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.945"
$ws.Range("E46").Value = "  -15.92%  "

# Row 47
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05402"
$ws.Range("E47").Value = "  -6.24%  "

# Row 48
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1114"
$ws.Range("E48").Value = "  -10.70%  "

# Row 49
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.036"
$ws.Range("E49").Value = "  -18.47%  "

# Row 50
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "30.25"
$ws.Range("E50").Value = "  -12.95%  "

# Row 51
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "52.93"
$ws.Range("E51").Value = "  -11.83%  "
